$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting (values such as
# "26.787.29" or "216.82" would otherwise be auto-coerced to numbers).
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '26.787.29'
$ws.Range('E2').Value = '  +0.07%  '

$ws.Range('D3').Value = '1.637.99'
$ws.Range('E3').Value = '  +0.11%  '

$ws.Range('E4').Value = '  -0.53%  '

$ws.Range('D5').Value = '216.82'
$ws.Range('E5').Value = '  -1.25%  '

$ws.Range('D6').Value = '0.508'
$ws.Range('E6').Value = '  +2.12%  '

$ws.Range('E7').Value = '  -0.53%  '

$ws.Range('D8').Value = '0.254'
$ws.Range('E8').Value = '  +1.78%  '

$ws.Range('E9').Value = '  +0.27%  '

$ws.Range('D10').Value = '19.85'
$ws.Range('E10').Value = '  +4.03%  '

$ws.Range('D11').Value = '0.0845'
$ws.Range('E11').Value = '  +0.08%  '

$ws.Range('D12').Value = '1.866.51'
$ws.Range('E12').Value = '  +0.31%  '

$ws.Range('D13').Value = '1.639.11'
$ws.Range('E13').Value = '  +0.55%  '

$ws.Range('E14').Value = '  -0.10%  '

$ws.Range('D15').Value = '0.529'
$ws.Range('E15').Value = '  +1.07%  '

$ws.Range('D16').Value = '66.33'
$ws.Range('E16').Value = '  +2.95%  '

$ws.Range('D17').Value = '26.793.06'
$ws.Range('E17').Value = '  +0.25%  '

$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('E18').Value = '  -0.55%  '

$ws.Range('D19').Value = '218.76'
$ws.Range('E19').Value = '  +1.79%  '

$ws.Range('E20').Value = '  -0.64%  '

$ws.Range('D21').Value = '6.67'
$ws.Range('E21').Value = '  +6.78%  '

$ws.Range('E22').Value = '  +1.00%  '

$ws.Range('D23').Value = '2.44'
$ws.Range('E23').Value = '  +5.47%  '

$ws.Range('E24').Value = '  +0.70%  '

$ws.Range('D25').Value = '147.15'
$ws.Range('E25').Value = '  -0.45%  '

$ws.Range('E26').Value = '  -0.82%  '

$ws.Range('D27').Value = '7.37'
$ws.Range('E27').Value = '  +4.69%  '

$ws.Range('E28').Value = '  +0.45%  '

$ws.Range('E29').Value = '  +0.49%  '

$ws.Range('D30').Value = '0.0503'
$ws.Range('E30').Value = '  -0.10%  '

$ws.Range('E31').Value = '  -1.93%  '

$ws.Range('D32').Value = '3.33'
$ws.Range('E32').Value = '  -1.60%  '

$ws.Range('E33').Value = '  +0.98%  '

$ws.Range('D34').Value = '1.56'
$ws.Range('E34').Value = '  +1.70%  '

$ws.Range('D35').Value = '1.260.09'
$ws.Range('E35').Value = '  +0.30%  '

$ws.Range('E36').Value = '  -0.33%  '

$ws.Range('D37').Value = '0.0177'
$ws.Range('E37').Value = '  +1.06%  '

$ws.Range('E38').Value = '  +1.55%  '

$ws.Range('D39').Value = '0.834'
$ws.Range('E39').Value = '  +3.24%  '

$ws.Range('E40').Value = '  -0.70%  '

$ws.Range('D41').Value = '0.808'
$ws.Range('E41').Value = '  +0.69%  '

$ws.Range('D42').Value = '5.45'
$ws.Range('E42').Value = '  +3.03%  '

$ws.Range('D43').Value = '1.776.78'
$ws.Range('E43').Value = '  +0.36%  '

$ws.Range('D44').Value = '61.69'
$ws.Range('E44').Value = '  +3.25%  '

$ws.Range('D45').Value = '2.10'
$ws.Range('E45').Value = '  -0.69%  '

$ws.Range('D46').Value = '91.65'
$ws.Range('E46').Value = '  -0.32%  '

$ws.Range('D47').Value = '1.57'
$ws.Range('E47').Value = '  -1.06%  '

$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('E48').Value = '  +1.45%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0513'
$ws.Range('E49').Value = '  -0.61%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.62'
$ws.Range('E50').Value = '  +2.27%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.0963'
$ws.Range('E51').Value = '  +0.52%  '
